# Apply the "Diagramme avec References" update:
#  - Sheet "Metadata": bump the Date value.
#  - Sheet "Elements": add a new row describing the
#    ActiviteSoumiseReconnaissance.EntiteGeographique reference element, and
#    widen column K (Type(s)) so the new long URL value fits.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Metadata sheet: refresh the generation Date value.
# ---------------------------------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B8").Value = "2025-10-30T16:36:55+00:00"

# ---------------------------------------------------------------------------
# 2) Elements sheet: append row 12 describing the new reference element.
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Elements")

# Text-like values that would otherwise be auto-coerced to numbers by the
# `.Value` setter (mirrors typing a leading apostrophe in Excel's UI to force
# text entry) so they land as shared-string cells, matching the other rows.
$ws.Range("A12").Value = "ActiviteSoumiseReconnaissance.EntiteGeographique"
$ws.Range("B12").Value = "ActiviteSoumiseReconnaissance.EntiteGeographique"
$ws.Range("D12").Value = "'"
$ws.Range("F12").Value = "'1"
$ws.Range("G12").Value = "'1"
$ws.Range("H12").Value = "'"
$ws.Range("I12").Value = "'"
$ws.Range("J12").Value = "'"
$ws.Range("K12").Value = "https://interop.esante.gouv.fr/ig/mos/StructureDefinition/EntiteGeographique`n"
$ws.Range("L12").Value = "Lien vers la classe EntiteGeographique"
$ws.Range("M12").Value = "Lien vers la classe EntiteGeographique"
$ws.Range("P12").Value = "'"
$ws.Range("R12").Value = "'"
$ws.Range("S12").Value = "'"
$ws.Range("T12").Value = "'"
$ws.Range("U12").Value = "'"
$ws.Range("V12").Value = "'"
$ws.Range("W12").Value = "'"
$ws.Range("X12").Value = "'"
$ws.Range("Y12").Value = "'"
$ws.Range("Z12").Value = "'"
$ws.Range("AA12").Value = "'"
$ws.Range("AB12").Value = "'"
$ws.Range("AC12").Value = "'"
$ws.Range("AD12").Value = "'"
$ws.Range("AE12").Value = "'"
$ws.Range("AF12").Value = "ActiviteSoumiseReconnaissance.EntiteGeographique"
$ws.Range("AG12").Value = "'1"
$ws.Range("AH12").Value = "'1"
$ws.Range("AI12").Value = "'"
$ws.Range("AJ12").Value = "'"

# Columns C, E, N, O and Q stay completely blank in row 12, same as row 11.

# Re-apply row 11's cell formatting onto row 12 so every cell shares the same
# style (border/fill/wrap) as the rest of the table instead of the default
# style picked up above.
$ws.Range("A11:AJ11").Copy()
$ws.Range("A12:AJ12").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Widen column K ("Type(s)") to fit the new long reference URL.
$ws.Columns.Item(11).ColumnWidth = 61.65625

Write-Host "Row 12 (EntiteGeographique) added and Date refreshed."
